$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 52-61 in column F ("type") were left blank; fill them with "c"
# (character type), matching the pattern used by every other row in the
# "Formats RAPSS ANO 17" layout sheet.
for ($r = 52; $r -le 61; $r++) {
    $ws.Cells.Item($r, 6).Value = "c"
}
